# Updates cryptos list with latest scraped prices/volume percentages.
# Mirrors the commit "Updated cryptos list ... with GitHub Actions".
#
# Price cells (column D) must stay stored as text (they use "." as a
# thousands separator, e.g. "44.239.91"), so for every D write we force
# text storage via NumberFormat "@" and then restore the default style
# so we don't leave a lingering custom number format on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "44.239.91"
$ws.Range("E2").Value = "  +0.96%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.242.72"
$ws.Range("E3").Value = "  +0.33%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "305.42"
$ws.Range("E5").Value = "  -2.92%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "95.40"
$ws.Range("E6").Value = "  -2.68%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -0.15%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.523"
$ws.Range("E9").Value = "  -1.29%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "34.77"
$ws.Range("E10").Value = "  -2.62%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.17%  "

# Row 12 - Polkadot
Set-TextValue $ws.Range("D12") "7.21"
$ws.Range("E12").Value = "  -1.56%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "2.584.06"
$ws.Range("E14").Value = "  +0.36%  "

# Row 15 - WrappedEther
Set-TextValue $ws.Range("D15") "2.328.05"
$ws.Range("E15").Value = "  +3.93%  "

# Row 16 - Polygon
Set-TextValue $ws.Range("D16") "0.830"
$ws.Range("E16").Value = "  -0.45%  "

# Row 17 - Chainlink
Set-TextValue $ws.Range("D17") "13.55"
$ws.Range("E17").Value = "  -2.37%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "44.007.43"
$ws.Range("E18").Value = "  +0.77%  "

# Row 19 - ShibaInu
Set-TextValue $ws.Range("D19") "0.0₃0962"
$ws.Range("E19").Value = "  -0.24%  "

# Row 20 - Uniswap
Set-TextValue $ws.Range("D20") "6.36"
$ws.Range("E20").Value = "  +1.34%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "12.04"
$ws.Range("E21").Value = "  -7.77%  "

# Row 22 - Litecoin
Set-TextValue $ws.Range("D22") "65.51"
$ws.Range("E22").Value = "  -0.59%  "

# Row 23 - was PancakeSwap, becomes BitcoinCash
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D23") "237.77"
$ws.Range("E23").Value = "  +0.76%  "

# Row 24 - was BitcoinCash, becomes PancakeSwap
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D24") "2.94"
$ws.Range("E24").Value = "  -1.40%  "

# Row 25 - ImmutableX
Set-TextValue $ws.Range("D25") "1.99"
$ws.Range("E25").Value = "  -1.25%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.03%  "

# Row 27 - was InjectiveProtocol, becomes Cosmos
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "9.90"
$ws.Range("E27").Value = "  -1.53%  "

# Row 28 - was Cosmos, becomes Toncoin
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D28") "2.20"
$ws.Range("E28").Value = "  +2.49%  "

# Row 29 - was Toncoin, becomes InjectiveProtocol
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D29") "38.02"
$ws.Range("E29").Value = "  +4.07%  "

# Row 30 - EthereumClassic
Set-TextValue $ws.Range("D30") "20.04"
$ws.Range("E30").Value = "  +0.30%  "

# Row 31 - Filecoin
Set-TextValue $ws.Range("D31") "5.86"
$ws.Range("E31").Value = "  -1.56%  "

# Row 32 - Monero
Set-TextValue $ws.Range("D32") "152.54"
$ws.Range("E32").Value = "  -2.27%  "

# Row 33 - Hedera
Set-TextValue $ws.Range("D33") "0.0792"
$ws.Range("E33").Value = "  -4.25%  "

# Row 34 - WEMIXToken
Set-TextValue $ws.Range("D34") "2.62"
$ws.Range("E34").Value = "  -0.70%  "

# Row 35 - LidoDAOToken
Set-TextValue $ws.Range("D35") "3.21"
$ws.Range("E35").Value = "  -3.75%  "

# Row 36 - Stellar
Set-TextValue $ws.Range("D36") "0.119"
$ws.Range("E36").Value = "  +2.00%  "

# Row 37 - Kaspa
$ws.Range("E37").Value = "  -2.16%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  -7.50%  "

# Row 39 - NEARProtocol
Set-TextValue $ws.Range("D39") "3.55"
$ws.Range("E39").Value = "  +1.15%  "

# Row 40 - RenderToken
Set-TextValue $ws.Range("D40") "3.83"
$ws.Range("E40").Value = "  -3.76%  "

# Row 41 - Celestia
Set-TextValue $ws.Range("D41") "14.27"
$ws.Range("E41").Value = "  -7.97%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  -2.48%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.20%  "

# Row 44 - Maker
Set-TextValue $ws.Range("D44") "1.751.02"
$ws.Range("E44").Value = "  +3.04%  "

# Row 45 - BitcoinSV
Set-TextValue $ws.Range("D45") "82.52"
$ws.Range("E45").Value = "  +0.50%  "

# Row 46 - Algorand
$ws.Range("E46").Value = "  -1.68%  "

# Row 47 - Aave
Set-TextValue $ws.Range("D47") "99.75"
$ws.Range("E47").Value = "  -1.58%  "

# Row 48 - THORChain
Set-TextValue $ws.Range("D48") "4.93"
$ws.Range("E48").Value = "  -3.93%  "

# Row 49 - was Stacks, becomes FraxShare
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D49") "8.09"
$ws.Range("E49").Value = "  -0.42%  "

# Row 50 - was FraxShare, becomes Stacks
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D50") "1.58"
$ws.Range("E50").Value = "  -1.55%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "54.48"
$ws.Range("E51").Value = "  -2.81%  "
